{"js": "// Apply the per-cell text replacements described by the diff.\n// Each (old, new) pair is unique and appears exactly once in the\n// document, so a plain \"search then replace\" per pair is unambiguous.\nconst replacements = [\n  [\"2024-02-24 Saturday\", \"2024-02-25 Sunday\"],\n  [\"86\u00d746=\", \"28\u00d797=\"],\n  [\"49\u00d730=\", \"40\u00d783=\"],\n  [\"46\u00d776=\", \"66\u00d774=\"],\n  [\"25\u00d744=\", \"89\u00d725=\"],\n  [\"65\u00d761=\", \"89\u00d748=\"],\n  [\"76\u00d786=\", \"81\u00d775=\"],\n  [\"57\u00d734=\", \"65\u00d726=\"],\n  [\"15\u00d760=\", \"40\u00d762=\"],\n  [\"87\u00d714=\", \"26\u00d753=\"],\n  [\"72\u00d762=\", \"70\u00d765=\"],\n  [\"15\u00d788=\", \"92\u00d799=\"],\n  [\"22\u00d755=\", \"36\u00d773=\"],\n  [\"23\u00d785=\", \"61\u00d755=\"],\n  [\"13\u00d760=\", \"96\u00d792=\"],\n  [\"85\u00d734=\", \"32\u00d734=\"],\n  [\"25\u00d794=\", \"46\u00d773=\"],\n  [\"92\u00d724=\", \"49\u00d788=\"],\n  [\"42\u00d764=\", \"78\u00d763=\"],\n  [\"11\u00d795=\", \"47\u00d756=\"],\n  [\"77\u00d799=\", \"26\u00d766=\"],\n  [\"48\u00d772=\", \"61\u00d744=\"],\n  [\"30\u00d754=\", \"61\u00d766=\"],\n  [\"86\u00d727=\", \"18\u00d783=\"],\n  [\"57\u00d766=\", \"60\u00d715=\"],\n  [\"83\u00d735=\", \"15\u00d790=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the per-cell text replacements described by the diff.\n# Each (old, new) pair is unique and appears exactly once in the\n# document, so Find/Replace across the whole document body is\n# unambiguous for every pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-24 Saturday\", \"2024-02-25 Sunday\"),\n    @(\"86\u00d746=\", \"28\u00d797=\"),\n    @(\"49\u00d730=\", \"40\u00d783=\"),\n    @(\"46\u00d776=\", \"66\u00d774=\"),\n    @(\"25\u00d744=\", \"89\u00d725=\"),\n    @(\"65\u00d761=\", \"89\u00d748=\"),\n    @(\"76\u00d786=\", \"81\u00d775=\"),\n    @(\"57\u00d734=\", \"65\u00d726=\"),\n    @(\"15\u00d760=\", \"40\u00d762=\"),\n    @(\"87\u00d714=\", \"26\u00d753=\"),\n    @(\"72\u00d762=\", \"70\u00d765=\"),\n    @(\"15\u00d788=\", \"92\u00d799=\"),\n    @(\"22\u00d755=\", \"36\u00d773=\"),\n    @(\"23\u00d785=\", \"61\u00d755=\"),\n    @(\"13\u00d760=\", \"96\u00d792=\"),\n    @(\"85\u00d734=\", \"32\u00d734=\"),\n    @(\"25\u00d794=\", \"46\u00d773=\"),\n    @(\"92\u00d724=\", \"49\u00d788=\"),\n    @(\"42\u00d764=\", \"78\u00d763=\"),\n    @(\"11\u00d795=\", \"47\u00d756=\"),\n    @(\"77\u00d799=\", \"26\u00d766=\"),\n    @(\"48\u00d772=\", \"61\u00d744=\"),\n    @(\"30\u00d754=\", \"61\u00d766=\"),\n    @(\"86\u00d727=\", \"18\u00d783=\"),\n    @(\"57\u00d766=\", \"60\u00d715=\"),\n    @(\"83\u00d735=\", \"15\u00d790=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
